# Auto-generated script applying scheduled market-price refresh to Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 15637681
$ws.Range("I11").Value = 15637681
$ws.Range("K11").Value = 15637681
$ws.Range("M11").Value = -15637541
$ws.Range("H86").Value = 951.25
$ws.Range("I86").Value = 904.17645
$ws.Range("J86").Value = 1218
$ws.Range("K86").Value = 904.17645
$ws.Range("L86").Value = 1218
$ws.Range("M86").Value = 218.82355
$ws.Range("N86").Value = -3464
$ws.Range("H89").Value = 951.25
$ws.Range("I89").Value = 904.17645
$ws.Range("J89").Value = 1218
$ws.Range("K89").Value = 4520.882250000001
$ws.Range("L89").Value = 6090
$ws.Range("M89").Value = 1095.117749999999
$ws.Range("N89").Value = -17322
$ws.Range("H141").Value = 7509.1665
$ws.Range("I141").Value = 2735
$ws.Range("J141").Value = 9896.25
$ws.Range("K141").Value = 8205
$ws.Range("L141").Value = 29688.75
$ws.Range("M141").Value = -3025
$ws.Range("N141").Value = -40048.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 10000
$ws.Range("J19").Value = 10000
$ws.Range("L19").Value = 10000
$ws.Range("N19").Value = -10458

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 33335622
$ws.Range("I105").Value = 50002070
$ws.Range("K105").Value = 50002070
$ws.Range("M105").Value = -50000323

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 58427.895
$ws.Range("I134").Value = 1506.5714
$ws.Range("J134").Value = 217807.6
$ws.Range("K134").Value = 4519.7142
$ws.Range("L134").Value = 653422.8
$ws.Range("M134").Value = -1984.7142
$ws.Range("N134").Value = -658492.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 784
$ws.Range("I5").Value = 306.66666
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 919.9999799999999
$ws.Range("L5").Value = 4500
$ws.Range("M5").Value = -807.9999799999999
$ws.Range("N5").Value = -4724
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16372
$ws.Range("H63").Value = 3391.5557
$ws.Range("I63").Value = 3127.5
$ws.Range("J63").Value = 3602.8
$ws.Range("K63").Value = 9382.5
$ws.Range("L63").Value = 10808.4
$ws.Range("M63").Value = -8633.5
$ws.Range("N63").Value = -12306.4
$ws.Range("H64").Value = 3304.08
$ws.Range("I64").Value = 1753
$ws.Range("J64").Value = 3599.524
$ws.Range("K64").Value = 5259
$ws.Range("L64").Value = 10798.572
$ws.Range("M64").Value = -4989
$ws.Range("N64").Value = -11338.572
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51864
$ws.Range("H66").Value = 3391.5557
$ws.Range("I66").Value = 3127.5
$ws.Range("J66").Value = 3602.8
$ws.Range("K66").Value = 28147.5
$ws.Range("L66").Value = 32425.2
$ws.Range("M66").Value = -24403.5
$ws.Range("N66").Value = -39913.2
$ws.Range("H67").Value = 3304.08
$ws.Range("I67").Value = 1753
$ws.Range("J67").Value = 3599.524
$ws.Range("K67").Value = 5259
$ws.Range("L67").Value = 10798.572
$ws.Range("M67").Value = -4323
$ws.Range("N67").Value = -12670.572
$ws.Range("H68").Value = 743.7012999999999
$ws.Range("I68").Value = 667.08826
$ws.Range("J68").Value = 1322.5555
$ws.Range("K68").Value = 2001.26478
$ws.Range("L68").Value = 3967.6665
$ws.Range("M68").Value = -1190.26478
$ws.Range("N68").Value = -5589.666499999999
$ws.Range("H69").Value = 799.75
$ws.Range("I69").Value = 799.75
$ws.Range("K69").Value = 2399.25
$ws.Range("M69").Value = -1588.25
$ws.Range("H71").Value = 743.7012999999999
$ws.Range("I71").Value = 667.08826
$ws.Range("J71").Value = 1322.5555
$ws.Range("K71").Value = 6003.79434
$ws.Range("L71").Value = 11902.9995
$ws.Range("M71").Value = -1947.79434
$ws.Range("N71").Value = -20014.9995
$ws.Range("H72").Value = 799.75
$ws.Range("I72").Value = 799.75
$ws.Range("K72").Value = 7197.75
$ws.Range("M72").Value = -3141.75
$ws.Range("H74").Value = 1980
$ws.Range("I74").Value = 1980
$ws.Range("K74").Value = 5940
$ws.Range("M74").Value = -4879
$ws.Range("H75").Value = 1124.3
$ws.Range("I75").Value = 1092.875
$ws.Range("J75").Value = 1250
$ws.Range("K75").Value = 3278.625
$ws.Range("L75").Value = 3750
$ws.Range("M75").Value = -2280.625
$ws.Range("N75").Value = -5746
$ws.Range("H77").Value = 1980
$ws.Range("I77").Value = 1980
$ws.Range("K77").Value = 17820
$ws.Range("M77").Value = -12516
$ws.Range("H78").Value = 1124.3
$ws.Range("I78").Value = 1092.875
$ws.Range("J78").Value = 1250
$ws.Range("K78").Value = 9835.875
$ws.Range("L78").Value = 11250
$ws.Range("M78").Value = -4843.875
$ws.Range("N78").Value = -21234
$ws.Range("H80").Value = 2121.3684
$ws.Range("I80").Value = 2400
$ws.Range("J80").Value = 2021.8572
$ws.Range("K80").Value = 7200
$ws.Range("L80").Value = 6065.571599999999
$ws.Range("M80").Value = -6264
$ws.Range("N80").Value = -7937.571599999999
$ws.Range("H81").Value = 11303.4
$ws.Range("I81").Value = 1009
$ws.Range("J81").Value = 26745
$ws.Range("K81").Value = 3027
$ws.Range("L81").Value = 80235
$ws.Range("M81").Value = -1904
$ws.Range("N81").Value = -82481
$ws.Range("H83").Value = 2121.3684
$ws.Range("I83").Value = 2400
$ws.Range("J83").Value = 2021.8572
$ws.Range("K83").Value = 21600
$ws.Range("L83").Value = 18196.7148
$ws.Range("M83").Value = -16920
$ws.Range("N83").Value = -27556.7148
$ws.Range("H84").Value = 11303.4
$ws.Range("I84").Value = 1009
$ws.Range("J84").Value = 26745
$ws.Range("K84").Value = 9081
$ws.Range("L84").Value = 240705
$ws.Range("M84").Value = -3465
$ws.Range("N84").Value = -251937
$ws.Range("H113").Value = 525.8095
$ws.Range("J113").Value = 577
$ws.Range("L113").Value = 1731
$ws.Range("N113").Value = -6071
$ws.Range("H135").Value = 784
$ws.Range("I135").Value = 306.66666
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 2759.99994
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -224.9999399999997
$ws.Range("N135").Value = -18570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4102.375
$ws.Range("I80").Value = 3760
$ws.Range("J80").Value = 4258
$ws.Range("K80").Value = 3760
$ws.Range("L80").Value = 4258
$ws.Range("M80").Value = -2762
$ws.Range("N80").Value = -6254
$ws.Range("H83").Value = 4102.375
$ws.Range("I83").Value = 3760
$ws.Range("J83").Value = 4258
$ws.Range("K83").Value = 18800
$ws.Range("L83").Value = 21290
$ws.Range("M83").Value = -13808
$ws.Range("N83").Value = -31274

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 3000
$ws.Range("J25").Value = 3000
$ws.Range("L25").Value = 3000
$ws.Range("N25").Value = -3460
$ws.Range("H69").Value = 40163
$ws.Range("J69").Value = 40163
$ws.Range("L69").Value = 40163
$ws.Range("N69").Value = -41785
$ws.Range("H72").Value = 40163
$ws.Range("J72").Value = 40163
$ws.Range("L72").Value = 120489
$ws.Range("N72").Value = -128601

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 70003.5
$ws.Range("J18").Value = 70003.5
$ws.Range("L18").Value = 70003.5
$ws.Range("N18").Value = -70349.5
$ws.Range("H132").Value = 78586.19500000001
$ws.Range("I132").Value = 54067.79
$ws.Range("J132").Value = 145136.14
$ws.Range("K132").Value = 162203.37
$ws.Range("L132").Value = 435408.42
$ws.Range("M132").Value = -159673.37
$ws.Range("N132").Value = -440468.42

Write-Output "Applied scheduled market-price refresh."